$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.796.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("E4").Value = "  -0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5029"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.47%  "

# Row 8
$ws.Range("E8").Value = "  -0.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07695"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "

# Row 12
$ws.Range("E12").Value = "  -0.68%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.857.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.31%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.601.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5426"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7925"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.824.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.32%  "

# Row 19
$ws.Range("E19").Value = "  -0.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.332"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.914"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.950"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "

# Row 24
$ws.Range("E24").Value = "  -0.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.931"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "

# Row 27
$ws.Range("E27").Value = "  -3.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.712"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.239"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.262"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.188"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.541"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.54%  "

# Row 35
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.172.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.626"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8918"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5583"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.697"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8070"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.770.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.28%  "

# Row 46
$ws.Range("E46").Value = "  -0.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4513"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("E49").Value = "  -1.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05076"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
